# #33 Removed Medium Prio. A task will have only 2 priorities: High and Low
#
# Task Priority table (column D, rows 4-7):
#   D4 = 00=High     (unchanged)
#   D5 = 01=Medium  -> 01=Low
#   D6 = 10=Low     -> 10=Not Used
#   D7 = 11=Not Used (unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "01=Low"
$ws.Range("D6").Value = "10=Not Used"

$ws.Range("D6").Select()
